$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.228.78"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "2.927.57"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "595.05"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "143.40"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -1.24%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("E10").Value = "  -2.14%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.436"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  -1.16%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "33.22"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "3.412.79"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").Value = "61.240.74"

$ws.Range("D17").Value = "2.930.24"
$ws.Range("E17").Value = "  +0.40%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.64"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "432.05"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.32%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.50"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.672"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.36%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.04"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "81.57"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.81"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.01%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.59%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.72"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("E28").Value = "  -4.36%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.61%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "26.61"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("E32").Value = "  +1.35%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "0.0₃0870"
$ws.Range("E34").Value = "  +2.23%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -0.36%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.55%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("E39").Value = "  -0.40%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.50"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "41.93"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.88%  "

$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("D44").Value = "2.691.06"
$ws.Range("E44").Value = "  -0.35%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "133.77"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "361.49"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.48"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.12%  "

$ws.Range("E49").Value = "  -1.25%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("E51").Value = "  -0.66%  "
